# Apply the cryptos-list refresh described in the commit:
#  - update Price/Volume(1h) figures for the existing 35 rows that keep their coin
#  - insert "Frax" as a new row 37, shifting HuobiToken..Maker down by one row,
#    which pushes "Cronos" off the bottom of the 50-row table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "30.429.56"
$ws.Range('E2').Value = "  +0.94%  "
$ws.Range('D3').Value = "1.870.75"
$ws.Range('E3').Value = "  +0.55%  "
$ws.Range('E4').Value = "  +0.13%  "
$ws.Range('D5').Value = "'246.75"
$ws.Range('E5').Value = "  +2.41%  "
$ws.Range('E6').Value = "  +0.11%  "
$ws.Range('D7').Value = "'0.4733"
$ws.Range('E7').Value = "  +0.83%  "
$ws.Range('D8').Value = "'0.2904"
$ws.Range('E8').Value = "  +1.63%  "
$ws.Range('D9').Value = "'0.06497"
$ws.Range('E9').Value = "  +0.43%  "
$ws.Range('D10').Value = "'21.92"
$ws.Range('E10').Value = "  +6.50%  "
$ws.Range('D11').Value = "'0.07719"
$ws.Range('E11').Value = "  +0.72%  "
$ws.Range('D12').Value = "'97.49"
$ws.Range('E12').Value = "  +4.05%  "
$ws.Range('D13').Value = "'0.7376"
$ws.Range('E13').Value = "  +8.58%  "
$ws.Range('D14').Value = "1.872.38"
$ws.Range('E14').Value = "  +0.66%  "
$ws.Range('D15').Value = "'5.111"
$ws.Range('E15').Value = "  +1.01%  "
$ws.Range('D16').Value = "'273.05"
$ws.Range('E16').Value = "  +1.81%  "
$ws.Range('D17').Value = "30.415.07"
$ws.Range('E17').Value = "  +0.94%  "
$ws.Range('D18').Value = "'13.38"
$ws.Range('E18').Value = "  +0.74%  "
$ws.Range('D19').Value = "'0.000007537"
$ws.Range('E19').Value = "  +0.19%  "
$ws.Range('E20').Value = "  +0.12%  "
$ws.Range('D21').Value = "2.126.39"
$ws.Range('E21').Value = "  +1.45%  "
$ws.Range('E22').Value = "  +0.20%  "
$ws.Range('D23').Value = "'5.220"
$ws.Range('E23').Value = "  +1.29%  "
$ws.Range('E24').Value = "  +1.18%  "
$ws.Range('E25').Value = "  -0.62%  "
$ws.Range('D26').Value = "'163.81"
$ws.Range('E26').Value = "  -1.37%  "
$ws.Range('D27').Value = "'18.82"
$ws.Range('E27').Value = "  +0.61%  "
$ws.Range('E28').Value = "  +2.83%  "
$ws.Range('E29').Value = "  +2.32%  "
$ws.Range('D30').Value = "'1.364"
$ws.Range('E30').Value = "  -0.64%  "
$ws.Range('D31').Value = "'1.507"
$ws.Range('E31').Value = "  +0.60%  "
$ws.Range('D32').Value = "'4.299"
$ws.Range('E32').Value = "  +2.13%  "
$ws.Range('D33').Value = "'4.121"
$ws.Range('E33').Value = "  +3.34%  "
$ws.Range('D34').Value = "'0.04804"
$ws.Range('E34').Value = "  +2.96%  "
$ws.Range('D35').Value = "'1.121"
$ws.Range('E35').Value = "  +1.23%  "
$ws.Range('D36').Value = "'0.6988"
$ws.Range('E36').Value = "  +2.20%  "
$ws.Range('B37').Value = "Frax"
$ws.Range('C37').Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range('D37').Value = "'1.000"
$ws.Range('E37').Value = "  +0.12%  "
$ws.Range('B38').Value = "HuobiToken"
$ws.Range('C38').Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D38').Value = "'2.715"
$ws.Range('E38').Value = "  +0.04%  "
$ws.Range('B39').Value = "VeChain"
$ws.Range('C39').Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('D39').Value = "'0.01857"
$ws.Range('E39').Value = "  +2.10%  "
$ws.Range('B40').Value = "MXToken"
$ws.Range('C40').Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D40').Value = "'2.750"
$ws.Range('E40').Value = "  +1.09%  "
$ws.Range('B41').Value = "FraxShare"
$ws.Range('C41').Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D41').Value = "'6.269"
$ws.Range('E41').Value = "  -1.05%  "
$ws.Range('B42').Value = "RenderToken"
$ws.Range('C42').Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D42').Value = "'1.970"
$ws.Range('E42').Value = "  +4.97%  "
$ws.Range('B43').Value = "Aave"
$ws.Range('C43').Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D43').Value = "'70.96"
$ws.Range('E43').Value = "  +1.38%  "
$ws.Range('B44').Value = "TheSandbox"
$ws.Range('C44').Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('D44').Value = "'0.4188"
$ws.Range('E44').Value = "  +3.64%  "
$ws.Range('B45').Value = "PaxDollar"
$ws.Range('C45').Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('D45').Value = "'1.001"
$ws.Range('E45').Value = "  +0.16%  "
$ws.Range('B46').Value = "TrustWalletToken"
$ws.Range('C46').Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D46').Value = "'0.8369"
$ws.Range('E46').Value = "  +0.55%  "
$ws.Range('B47').Value = "Quant"
$ws.Range('C47').Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('D47').Value = "'102.91"
$ws.Range('E47').Value = "  +1.12%  "
$ws.Range('B48').Value = "EnergySwap"
$ws.Range('C48').Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D48').Value = "'9.257"
$ws.Range('E48').Value = "  +0.45%  "
$ws.Range('B49').Value = "Aptos"
$ws.Range('C49').Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('D49').Value = "'7.033"
$ws.Range('E49').Value = "  +1.78%  "
$ws.Range('B50').Value = "Elrond"
$ws.Range('C50').Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range('D50').Value = "'35.45"
$ws.Range('E50').Value = "  +3.78%  "
$ws.Range('B51').Value = "Maker"
$ws.Range('C51').Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D51').Value = "'920.70"
$ws.Range('E51').Value = "  -0.52%  "
